$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = -1
$ws.Range("F7").Value = 5
$ws.Range("F9").Value = -10
$ws.Range("F10").Value = -5
$ws.Range("F12").Value = 0
$ws.Range("F21").Value = 5
$ws.Range("F28").Value = 2
$ws.Range("F29").Value = -1
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = -2
$ws.Range("F38").Value = -1
$ws.Range("F40").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("F45").Value = -1
$ws.Range("F46").Value = -6
$ws.Range("F47").Value = -9
$ws.Range("F50").Value = 0
